$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.306379666666667
$ws.Range("H2").Value = 3.919139
$ws.Range("I2").Value = 0.4034923136874173
$ws.Range("J2").Value = 0.4034923136874172
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 14.96835266666667
$ws.Range("N2").Value = 44.905058
$ws.Range("O2").Value = 0.1240053612000741
$ws.Range("P2").Value = 0.1240053612000741
$ws.Range("Q2").Value = 19.55435156722911
$ws.Range("R2").Value = 175.989164105062
$ws.Range("S2").Value = 0.05003521010026179
$ws.Range("T2").Value = 0.05003521010026177

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.306379666666667
$ws.Range("H3").Value = 3.919139
$ws.Range("I3").Value = 0.4034923136874173
$ws.Range("J3").Value = 0.4034923136874172
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 80.77474466666666
$ws.Range("N3").Value = 242.324234
$ws.Range("O3").Value = 0.6691786071115035
$ws.Range("P3").Value = 0.6691786071115035
$ws.Range("Q3").Value = 105.5224840127251
$ws.Range("R3").Value = 949.702356114526
$ws.Range("S3").Value = 0.2700084244535437
$ws.Range("T3").Value = 0.2700084244535437

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.306379666666667
$ws.Range("H4").Value = 3.919139
$ws.Range("I4").Value = 0.4034923136874173
$ws.Range("J4").Value = 0.4034923136874172
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 24.96420533333334
$ws.Range("N4").Value = 74.892616
$ws.Range("O4").Value = 0.2068160316884225
$ws.Range("P4").Value = 0.2068160316884225
$ws.Range("Q4").Value = 32.61273024195823
$ws.Range("R4").Value = 293.514572177624
$ws.Range("S4").Value = 0.0834486791336118
$ws.Range("T4").Value = 0.08344867913361179

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.072562
$ws.Range("H5").Value = 3.217686
$ws.Range("I5").Value = 0.3312746929515923
$ws.Range("J5").Value = 0.3312746929515923
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 14.96835266666667
$ws.Range("N5").Value = 44.905058
$ws.Range("O5").Value = 0.1240053612000741
$ws.Range("P5").Value = 0.1240053612000741
$ws.Range("Q5").Value = 16.05448627286533
$ws.Range("R5").Value = 144.490376455788
$ws.Range("S5").Value = 0.04107983795590585
$ws.Range("T5").Value = 0.04107983795590585

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.072562
$ws.Range("H6").Value = 3.217686
$ws.Range("I6").Value = 0.3312746929515923
$ws.Range("J6").Value = 0.3312746929515923
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 80.77474466666666
$ws.Range("N6").Value = 242.324234
$ws.Range("O6").Value = 0.6691786071115035
$ws.Range("P6").Value = 0.6691786071115035
$ws.Range("Q6").Value = 86.63592168916931
$ws.Range("R6").Value = 779.7232952025239
$ws.Range("S6").Value = 0.2216819376006376
$ws.Range("T6").Value = 0.2216819376006376

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.072562
$ws.Range("H7").Value = 3.217686
$ws.Range("I7").Value = 0.3312746929515923
$ws.Range("J7").Value = 0.3312746929515923
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 24.96420533333334
$ws.Range("N7").Value = 74.892616
$ws.Range("O7").Value = 0.2068160316884225
$ws.Range("P7").Value = 0.2068160316884225
$ws.Range("Q7").Value = 26.77565800073067
$ws.Range("R7").Value = 240.980922006576
$ws.Range("S7").Value = 0.06851291739504894
$ws.Range("T7").Value = 0.06851291739504896

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.8587400000000001
$ws.Range("H8").Value = 2.57622
$ws.Range("I8").Value = 0.2652329933609903
$ws.Range("J8").Value = 0.2652329933609903
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 14.96835266666667
$ws.Range("N8").Value = 44.905058
$ws.Range("O8").Value = 0.1240053612000741
$ws.Range("P8").Value = 0.1240053612000741
$ws.Range("Q8").Value = 12.85392316897333
$ws.Range("R8").Value = 115.68530852076
$ws.Range("S8").Value = 0.03289031314390646
$ws.Range("T8").Value = 0.03289031314390645

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.8587400000000001
$ws.Range("H9").Value = 2.57622
$ws.Range("I9").Value = 0.2652329933609903
$ws.Range("J9").Value = 0.2652329933609903
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 80.77474466666666
$ws.Range("N9").Value = 242.324234
$ws.Range("O9").Value = 0.6691786071115035
$ws.Range("P9").Value = 0.6691786071115035
$ws.Range("Q9").Value = 69.36450423505333
$ws.Range("R9").Value = 624.28053811548
$ws.Range("S9").Value = 0.1774882450573221
$ws.Range("T9").Value = 0.1774882450573221

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.8587400000000001
$ws.Range("H10").Value = 2.57622
$ws.Range("I10").Value = 0.2652329933609903
$ws.Range("J10").Value = 0.2652329933609903
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 24.96420533333334
$ws.Range("N10").Value = 74.892616
$ws.Range("O10").Value = 0.2068160316884225
$ws.Range("P10").Value = 0.2068160316884225
$ws.Range("Q10").Value = 21.43776168794667
$ws.Range("R10").Value = 192.93985519152
$ws.Range("S10").Value = 0.05485443515976172
$ws.Range("T10").Value = 0.05485443515976172

